# Updated cryptos list on Tue Oct 31 03:21:29 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Vol($row, $pct) {
    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = "  $pct  "
}

function Set-Price($row, $price) {
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $price
}

# Row 2 - Bitcoin
Set-Price 2 "34.425.72"
Set-Vol 2 "+0.19%"

# Row 3 - Ethereum
Set-Price 3 "1.805.42"
Set-Vol 3 "+1.05%"

# Row 4 - TetherUSD
Set-Vol 4 "+0.05%"

# Row 5 - BNB
Set-Price 5 "227.68"
Set-Vol 5 "+0.51%"

# Row 6 - XRP
Set-Price 6 "0.579"
Set-Vol 6 "+4.07%"

# Row 7 - USDC
Set-Vol 7 "+0.04%"

# Row 8 - Solana
Set-Price 8 "35.89"
Set-Vol 8 "+8.77%"

# Row 9 - Cardano
Set-Price 9 "0.302"
Set-Vol 9 "+2.25%"

# Row 10 - Dogecoin
Set-Price 10 "0.0694"
Set-Vol 10 "+0.80%"

# Row 11 - TRON
Set-Price 11 "0.0965"
Set-Vol 11 "+1.94%"

# Row 12 - WrappedliquidstakedEther2.0
Set-Price 12 "2.067.79"
Set-Vol 12 "+1.10%"

# Row 13 - Chainlink
Set-Price 13 "11.46"
Set-Vol 13 "+1.79%"

# Row 14 - WrappedEther
Set-Price 14 "1.811.25"
Set-Vol 14 "+1.22%"

# Row 15 - Polygon
Set-Price 15 "0.645"
Set-Vol 15 "+1.59%"

# Row 16 - Polkadot
Set-Price 16 "4.52"
Set-Vol 16 "+5.25%"

# Row 17 - WrappedBTC
Set-Price 17 "34.404.40"
Set-Vol 17 "+0.02%"

# Row 18 - Litecoin
Set-Price 18 "69.18"
Set-Vol 18 "+1.05%"

# Row 19 - BitcoinCash
Set-Price 19 "245.91"
Set-Vol 19 "+0.31%"

# Row 20 - ShibaInu
Set-Price 20 "0.0₃0797"
Set-Vol 20 "+0.22%"

# Row 21 - Avalanche
Set-Price 21 "11.49"
Set-Vol 21 "+1.63%"

# Row 22 - Dai
Set-Vol 22 "+0.00%"

# Row 23 - Uniswap
Set-Price 23 "4.20"
Set-Vol 23 "+0.85%"

# Row 24 - Toncoin
Set-Vol 24 "+3.21%"

# Row 25 - Stellar/etc.
Set-Price 25 "170.94"
Set-Vol 25 "+1.49%"

# Row 26
Set-Price 26 "7.91"
Set-Vol 26 "+7.73%"

# Row 27
Set-Price 27 "17.00"
Set-Vol 27 "+2.77%"

# Row 28
Set-Price 28 "0.119"
Set-Vol 28 "+3.02%"

# Row 29 - BinanceUSD
Set-Vol 29 "+0.03%"

# Row 30 - InternetComputer(DFINITY)
Set-Price 30 "4.06"
Set-Vol 30 "+0.68%"

# Row 31 - Hedera
Set-Price 31 "0.0532"
Set-Vol 31 "+1.11%"

# Row 32 - Filecoin
Set-Price 32 "3.86"
Set-Vol 32 "+1.44%"

# Row 33 - PancakeSwap
Set-Vol 33 "+0.86%"

# Row 34 - LidoDAOToken
Set-Price 34 "1.84"
Set-Vol 34 "+0.86%"

# Row 35 - Maker
Set-Price 35 "1.398.83"
Set-Vol 35 "-0.89%"

# Row 36 - ImmutableX
Set-Price 36 "0.674"
Set-Vol 36 "-1.64%"

# Row 37 - RenderToken
Set-Vol 37 "-3.20%"

# Row 38 - TrustWalletToken
Set-Vol 38 "-0.28%"

# Row 39 - VeChain
Set-Price 39 "0.0190"
Set-Vol 39 "-0.07%"

# Row 40 - WEMIXToken
Set-Price 40 "1.24"
Set-Vol 40 "+12.15%"

# Rows 41 & 42 swap: Aave <-> ARBITRUM
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-Price 41 "0.967"
Set-Vol 41 "+2.82%"

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-Price 42 "82.77"
Set-Vol 42 "-2.39%"

# Row 43 - MXToken
Set-Price 43 "2.82"
Set-Vol 43 "+1.85%"

# Row 44 - HuobiToken
Set-Vol 44 "+0.00%"

# Row 45 - InjectiveProtocol
Set-Price 45 "13.45"
Set-Vol 45 "-3.64%"

# Row 46 - FraxShare
Set-Price 46 "6.06"
Set-Vol 46 "-0.45%"

# Row 47 - Kaspa
Set-Vol 47 "-4.91%"

# Row 48 - RocketPoolETH
Set-Price 48 "1.967.06"
Set-Vol 48 "+1.07%"

# Row 49 - Quant
Set-Price 49 "105.29"
Set-Vol 49 "-0.10%"

# Row 50 - PaxDollar
Set-Vol 50 "+0.08%"

# Row 51 - BabyDogeCoin
Set-Vol 51 "+1.27%"
